$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 686.3125
$ws.Range("I9").Value = 335.83334
$ws.Range("J9").Value = 1737.75
$ws.Range("K9").Value = 335.83334
$ws.Range("L9").Value = 1737.75
$ws.Range("M9").Value = -166.83334
$ws.Range("N9").Value = -2075.75
$ws.Range("H18").Value = 1552.4706
$ws.Range("I18").Value = 399.57144
$ws.Range("K18").Value = 399.57144
$ws.Range("M18").Value = -115.57144
$ws.Range("H21").Value = 18500
$ws.Range("I21").Value = 18500
$ws.Range("K21").Value = 18500
$ws.Range("M21").Value = -18032
$ws.Range("H23").Value = 18500
$ws.Range("I23").Value = 18500
$ws.Range("K23").Value = 18500
$ws.Range("M23").Value = -18266
$ws.Range("H29").Value = 1560.6
$ws.Range("I29").Value = 334.33334
$ws.Range("J29").Value = 3400
$ws.Range("K29").Value = 1003.00002
$ws.Range("L29").Value = 10200
$ws.Range("M29").Value = -722.0000200000001
$ws.Range("N29").Value = -10762
$ws.Range("H33").Value = 2315.8125
$ws.Range("I33").Value = 465.5
$ws.Range("K33").Value = 465.5
$ws.Range("M33").Value = -236.5
$ws.Range("H38").Value = 3308.077
$ws.Range("I38").Value = 205.6
$ws.Range("J38").Value = 5247.125
$ws.Range("K38").Value = 616.8
$ws.Range("L38").Value = 15741.375
$ws.Range("M38").Value = -244.8
$ws.Range("N38").Value = -16485.375
$ws.Range("H53").Value = 627.1429000000001
$ws.Range("I53").Value = 315
$ws.Range("J53").Value = 2500
$ws.Range("K53").Value = 315
$ws.Range("L53").Value = 2500
$ws.Range("M53").Value = 322
$ws.Range("N53").Value = -3774
$ws.Range("H58").Value = 371
$ws.Range("I58").Value = 371
$ws.Range("K58").Value = 1113
$ws.Range("M58").Value = -963
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H62").Value = 4104.625
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 4104.625
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
$ws.Range("H76").Value = 6566.5
$ws.Range("J76").Value = 7199.875
$ws.Range("L76").Value = 7199.875
$ws.Range("N76").Value = -7829.875
$ws.Range("H79").Value = 6566.5
$ws.Range("J79").Value = 7199.875
$ws.Range("L79").Value = 7199.875
$ws.Range("N79").Value = -9383.875
$ws.Range("H87").Value = 67549.75
$ws.Range("J87").Value = 71057
$ws.Range("L87").Value = 71057
$ws.Range("N87").Value = -73553
$ws.Range("H90").Value = 67549.75
$ws.Range("J90").Value = 71057
$ws.Range("L90").Value = 213171
$ws.Range("N90").Value = -225651
$ws.Range("H111").Value = 75078.8
$ws.Range("J111").Value = 124665.555
$ws.Range("L111").Value = 373996.665
$ws.Range("N111").Value = -380130.665
$ws.Range("H112").Value = 6600.778
$ws.Range("J112").Value = 8031
$ws.Range("L112").Value = 24093
$ws.Range("N112").Value = -26309
$ws.Range("H125").Value = 126690.25
$ws.Range("I125").Value = 2348.5
$ws.Range("K125").Value = 21136.5
$ws.Range("M125").Value = -18676.5
$ws.Range("H136").Value = 94996.2
$ws.Range("J136").Value = 94996.2
$ws.Range("L136").Value = 94996.2
$ws.Range("N136").Value = -105196.2
$ws.Range("H137").Value = 2105.3784
$ws.Range("I137").Value = 1818.2122
$ws.Range("J137").Value = 4474.5
$ws.Range("K137").Value = 5454.6366
$ws.Range("L137").Value = 13423.5
$ws.Range("M137").Value = -2904.6366
$ws.Range("N137").Value = -18523.5
$ws.Range("H141").Value = 4268.75
$ws.Range("I141").Value = 3598.5
$ws.Range("J141").Value = 5385.8335
$ws.Range("K141").Value = 10795.5
$ws.Range("L141").Value = 16157.5005
$ws.Range("M141").Value = -5615.5
$ws.Range("N141").Value = -26517.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10753.151
$ws.Range("I32").Value = 4383.963
$ws.Range("K32").Value = 4383.963
$ws.Range("M32").Value = -4096.963
$ws.Range("H34").Value = 73500
$ws.Range("J34").Value = 87000
$ws.Range("L34").Value = 87000
$ws.Range("N34").Value = -87542
$ws.Range("H39").Value = 38907.2
$ws.Range("I39").Value = 6508
$ws.Range("J39").Value = 60506.668
$ws.Range("K39").Value = 6508
$ws.Range("L39").Value = 60506.668
$ws.Range("M39").Value = -5988
$ws.Range("N39").Value = -61546.668
$ws.Range("H45").Value = 1849.2858
$ws.Range("I45").Value = 1721.5
$ws.Range("J45").Value = 1900.4
$ws.Range("K45").Value = 1721.5
$ws.Range("L45").Value = 1900.4
$ws.Range("M45").Value = -1344.5
$ws.Range("N45").Value = -2654.4
$ws.Range("H74").Value = 5419.524
$ws.Range("I74").Value = 1855.8334
$ws.Range("J74").Value = 10171.111
$ws.Range("K74").Value = 1855.8334
$ws.Range("L74").Value = 10171.111
$ws.Range("M74").Value = -981.8334
$ws.Range("N74").Value = -11919.111
$ws.Range("H77").Value = 5419.524
$ws.Range("I77").Value = 1855.8334
$ws.Range("J77").Value = 10171.111
$ws.Range("K77").Value = 9279.166999999999
$ws.Range("L77").Value = 50855.55500000001
$ws.Range("M77").Value = -4911.166999999999
$ws.Range("N77").Value = -59591.55500000001
$ws.Range("H88").Value = 2901.9473
$ws.Range("I88").Value = 2157.2
$ws.Range("J88").Value = 3729.4443
$ws.Range("K88").Value = 2157.2
$ws.Range("L88").Value = 3729.4443
$ws.Range("M88").Value = -1751.2
$ws.Range("N88").Value = -4541.4443
$ws.Range("H91").Value = 2901.9473
$ws.Range("I91").Value = 2157.2
$ws.Range("J91").Value = 3729.4443
$ws.Range("K91").Value = 2157.2
$ws.Range("L91").Value = 3729.4443
$ws.Range("M91").Value = -753.1999999999998
$ws.Range("N91").Value = -6537.4443
$ws.Range("H122").Value = 2837.6553
$ws.Range("I122").Value = 1382.3043
$ws.Range("K122").Value = 4146.9129
$ws.Range("M122").Value = -1696.9129
$ws.Range("H132").Value = 4984.28
$ws.Range("I132").Value = 2128.2778
$ws.Range("K132").Value = 6384.8334
$ws.Range("M132").Value = -3854.8334
$ws.Range("H133").Value = 48086.11
$ws.Range("J133").Value = 48086.11
$ws.Range("L133").Value = 48086.11
$ws.Range("N133").Value = -53146.11
$ws.Range("H135").Value = 44637.5
$ws.Range("J135").Value = 44637.5
$ws.Range("L135").Value = 44637.5
$ws.Range("N135").Value = -54777.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 100000
$ws.Range("I59").Value = 100000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 100000
$ws.Range("L59").ClearContents()
$ws.Range("M59").Value = -99153
$ws.Range("N59").Value = 0
$ws.Range("H86").Value = 2959.48
$ws.Range("I86").Value = 946.7778
$ws.Range("J86").Value = 8135
$ws.Range("K86").Value = 946.7778
$ws.Range("L86").Value = 8135
$ws.Range("M86").Value = 176.2222
$ws.Range("N86").Value = -10381
$ws.Range("H89").Value = 2959.48
$ws.Range("I89").Value = 946.7778
$ws.Range("J89").Value = 8135
$ws.Range("K89").Value = 4733.889
$ws.Range("L89").Value = 40675
$ws.Range("M89").Value = 882.1109999999999
$ws.Range("N89").Value = -51907
$ws.Range("H99").Value = 2173.2
$ws.Range("I99").Value = 2025.7778
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2025.7778
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -527.7778000000001
$ws.Range("N99").Value = -6496
$ws.Range("H134").Value = 5575.6875
$ws.Range("I134").Value = 2363.6365
$ws.Range("J134").Value = 12642.2
$ws.Range("K134").Value = 7090.9095
$ws.Range("L134").Value = 37926.60000000001
$ws.Range("M134").Value = -4555.9095
$ws.Range("N134").Value = -42996.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2107399.8
$ws.Range("I6").Value = 3503333.2
$ws.Range("J6").Value = 13499.5
$ws.Range("K6").Value = 3503333.2
$ws.Range("L6").Value = 13499.5
$ws.Range("M6").Value = -3503220.2
$ws.Range("N6").Value = -13725.5
$ws.Range("H16").Value = 2240.5
$ws.Range("I16").Value = 700
$ws.Range("J16").Value = 2754
$ws.Range("K16").Value = 700
$ws.Range("L16").Value = 2754
$ws.Range("M16").Value = -413
$ws.Range("N16").Value = -3328
$ws.Range("H31").Value = 8895.125
$ws.Range("I31").Value = 4114.154
$ws.Range("J31").Value = 14545.363
$ws.Range("K31").Value = 4114.154
$ws.Range("L31").Value = 14545.363
$ws.Range("M31").Value = -3819.154
$ws.Range("N31").Value = -15135.363
$ws.Range("H34").Value = 8895.125
$ws.Range("I34").Value = 4114.154
$ws.Range("J34").Value = 14545.363
$ws.Range("K34").Value = 4114.154
$ws.Range("L34").Value = 14545.363
$ws.Range("M34").Value = -3912.154
$ws.Range("N34").Value = -14949.363
$ws.Range("H58").Value = 4645.2104
$ws.Range("I58").Value = 3475
$ws.Range("K58").Value = 3475
$ws.Range("M58").Value = -3272
$ws.Range("H62").Value = 3683.1765
$ws.Range("I62").Value = 3247.111
$ws.Range("J62").Value = 4173.75
$ws.Range("K62").Value = 3247.111
$ws.Range("L62").Value = 4173.75
$ws.Range("M62").Value = -2623.111
$ws.Range("N62").Value = -5421.75
$ws.Range("H65").Value = 3683.1765
$ws.Range("I65").Value = 3247.111
$ws.Range("J65").Value = 4173.75
$ws.Range("K65").Value = 16235.555
$ws.Range("L65").Value = 20868.75
$ws.Range("M65").Value = -13115.555
$ws.Range("N65").Value = -27108.75
$ws.Range("H99").Value = 2476.842
$ws.Range("I99").Value = 2476.842
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2476.842
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -978.8420000000001
$ws.Range("H113").Value = 2240.5
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 2754
$ws.Range("K113").Value = 700
$ws.Range("L113").Value = 2754
$ws.Range("M113").Value = 1470
$ws.Range("N113").Value = -7094
$ws.Range("H126").Value = 2476.842
$ws.Range("I126").Value = 2476.842
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7430.526
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -4960.526
$ws.Range("H132").Value = 6227.1816
$ws.Range("I132").Value = 4187.375
$ws.Range("J132").Value = 11666.667
$ws.Range("K132").Value = 12562.125
$ws.Range("L132").Value = 35000.001
$ws.Range("M132").Value = -10032.125
$ws.Range("N132").Value = -40060.001
$ws.Range("H134").Value = 6964.2354
$ws.Range("I134").Value = 5693
$ws.Range("J134").Value = 16498.5
$ws.Range("K134").Value = 17079
$ws.Range("L134").Value = 49495.5
$ws.Range("M134").Value = -14544
$ws.Range("N134").Value = -54565.5
$ws.Range("H135").Value = 69994
$ws.Range("J135").Value = 69994
$ws.Range("L135").Value = 69994
$ws.Range("N135").Value = -80134
$ws.Range("H136").Value = 4645.2104
$ws.Range("I136").Value = 3475
$ws.Range("K136").Value = 10425
$ws.Range("M136").Value = -7875
$ws.Range("H140").Value = 99196
$ws.Range("J140").Value = 99196
$ws.Range("L140").Value = 99196
$ws.Range("N140").Value = -109556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 804.25
$ws.Range("I2").Value = 116.36364
$ws.Range("J2").Value = 2317.6
$ws.Range("K2").Value = 698.18184
$ws.Range("L2").Value = 13905.6
$ws.Range("M2").Value = -585.18184
$ws.Range("N2").Value = -14131.6
$ws.Range("H68").Value = 1979
$ws.Range("I68").Value = 1968.125
$ws.Range("K68").Value = 5904.375
$ws.Range("M68").Value = -5093.375
$ws.Range("H71").Value = 1979
$ws.Range("I71").Value = 1968.125
$ws.Range("K71").Value = 17713.125
$ws.Range("M71").Value = -13657.125
$ws.Range("H97").Value = 740.75
$ws.Range("I97").Value = 683.3333
$ws.Range("J97").Value = 775.2
$ws.Range("K97").Value = 2049.9999
$ws.Range("L97").Value = 2325.6
$ws.Range("M97").Value = -1553.9999
$ws.Range("N97").Value = -3317.6
$ws.Range("H107").Value = 513.53845
$ws.Range("I107").Value = 202.66667
$ws.Range("J107").Value = 606.8
$ws.Range("K107").Value = 608.00001
$ws.Range("L107").Value = 1820.4
$ws.Range("M107").Value = 1311.99999
$ws.Range("N107").Value = -5660.4
$ws.Range("H137").Value = 3219.0256
$ws.Range("I137").Value = 2213.3572
$ws.Range("J137").Value = 3782.2
$ws.Range("K137").Value = 6640.071599999999
$ws.Range("L137").Value = 11346.6
$ws.Range("M137").Value = -1540.071599999999
$ws.Range("N137").Value = -21546.6
$ws.Range("H140").Value = 1718.8451
$ws.Range("I140").Value = 1147.25
$ws.Range("K140").Value = 3441.75
$ws.Range("M140").Value = 1738.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0
$ws.Range("H70").Value = 6215.6294
$ws.Range("I70").Value = 5632.3335
$ws.Range("J70").Value = 6507.278
$ws.Range("K70").Value = 5632.3335
$ws.Range("L70").Value = 6507.278
$ws.Range("M70").Value = -5362.3335
$ws.Range("N70").Value = -7047.278
$ws.Range("H73").Value = 6215.6294
$ws.Range("I73").Value = 5632.3335
$ws.Range("J73").Value = 6507.278
$ws.Range("K73").Value = 5632.3335
$ws.Range("L73").Value = 6507.278
$ws.Range("M73").Value = -4696.3335
$ws.Range("N73").Value = -8379.278
$ws.Range("H80").Value = 18249.625
$ws.Range("I80").Value = 17088.889
$ws.Range("J80").Value = 19742
$ws.Range("K80").Value = 17088.889
$ws.Range("L80").Value = 19742
$ws.Range("M80").Value = -16090.889
$ws.Range("N80").Value = -21738
$ws.Range("H83").Value = 18249.625
$ws.Range("I83").Value = 17088.889
$ws.Range("J83").Value = 19742
$ws.Range("K83").Value = 85444.44499999999
$ws.Range("L83").Value = 98710
$ws.Range("M83").Value = -80452.44499999999
$ws.Range("N83").Value = -108694
$ws.Range("H113").Value = 1987.2858
$ws.Range("I113").Value = 2430.3333
$ws.Range("J113").Value = 1655
$ws.Range("K113").Value = 2430.3333
$ws.Range("L113").Value = 1655
$ws.Range("M113").Value = -260.3332999999998
$ws.Range("N113").Value = -5995
$ws.Range("H122").Value = 5202.4443
$ws.Range("I122").Value = 1012.25
$ws.Range("J122").Value = 13582.833
$ws.Range("K122").Value = 3036.75
$ws.Range("L122").Value = 40748.499
$ws.Range("M122").Value = -586.75
$ws.Range("N122").Value = -45648.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 0
$ws.Range("H22").Value = 2336.3333
$ws.Range("I22").Value = 1004.75
$ws.Range("K22").Value = 1004.75
$ws.Range("M22").Value = -709.75
$ws.Range("H27").Value = 2336.3333
$ws.Range("I27").Value = 1004.75
$ws.Range("K27").Value = 1004.75
$ws.Range("M27").Value = -897.75
$ws.Range("H40").Value = 4635.769
$ws.Range("I40").Value = 3156.5
$ws.Range("K40").Value = 3156.5
$ws.Range("M40").Value = -3020.5
$ws.Range("H46").Value = 2789.2
$ws.Range("I46").Value = 639.5
$ws.Range("K46").Value = 639.5
$ws.Range("M46").Value = -451.5
$ws.Range("H68").Value = 2884.9
$ws.Range("I68").Value = 2860
$ws.Range("J68").Value = 2984.5
$ws.Range("K68").Value = 2860
$ws.Range("L68").Value = 2984.5
$ws.Range("M68").Value = -2111
$ws.Range("N68").Value = -4482.5
$ws.Range("H71").Value = 2884.9
$ws.Range("I71").Value = 2860
$ws.Range("J71").Value = 2984.5
$ws.Range("K71").Value = 14300
$ws.Range("L71").Value = 14922.5
$ws.Range("M71").Value = -10556
$ws.Range("N71").Value = -22410.5
$ws.Range("H82").Value = 1768.75
$ws.Range("I82").Value = 1287.75
$ws.Range("J82").Value = 2249.75
$ws.Range("K82").Value = 1287.75
$ws.Range("L82").Value = 2249.75
$ws.Range("M82").Value = -926.75
$ws.Range("N82").Value = -2971.75
$ws.Range("H85").Value = 1768.75
$ws.Range("I85").Value = 1287.75
$ws.Range("J85").Value = 2249.75
$ws.Range("K85").Value = 1287.75
$ws.Range("L85").Value = 2249.75
$ws.Range("M85").Value = -39.75
$ws.Range("N85").Value = -4745.75
$ws.Range("H93").Value = 2768.5454
$ws.Range("I93").Value = 2973.875
$ws.Range("J93").Value = 2221
$ws.Range("K93").Value = 2973.875
$ws.Range("L93").Value = 2221
$ws.Range("M93").Value = -1725.875
$ws.Range("N93").Value = -4717
$ws.Range("H94").Value = 58249.75
$ws.Range("I94").Value = 39999
$ws.Range("J94").Value = 64333.332
$ws.Range("K94").Value = 39999
$ws.Range("L94").Value = 64333.332
$ws.Range("M94").Value = -39323
$ws.Range("N94").Value = -65685.33199999999
$ws.Range("H100").Value = 4666.6665
$ws.Range("I100").Value = 5500
$ws.Range("K100").Value = 5500
$ws.Range("M100").Value = -4959
$ws.Range("H132").Value = 5519
$ws.Range("I132").Value = 3604
$ws.Range("J132").Value = 6795.6665
$ws.Range("K132").Value = 10812
$ws.Range("L132").Value = 20386.9995
$ws.Range("M132").Value = -8282
$ws.Range("N132").Value = -25446.9995
$ws.Range("H133").Value = 94942.71000000001
$ws.Range("J133").Value = 94942.71000000001
$ws.Range("L133").Value = 94942.71000000001
$ws.Range("N133").Value = -100002.71
$ws.Range("H136").Value = 6500.033
$ws.Range("I136").Value = 2946.9092
$ws.Range("J136").Value = 8557.105
$ws.Range("K136").Value = 8840.7276
$ws.Range("L136").Value = 25671.315
$ws.Range("M136").Value = -6290.7276
$ws.Range("N136").Value = -30771.315
$ws.Range("H140").Value = 60628
$ws.Range("J140").Value = 60628
$ws.Range("L140").Value = 60628
$ws.Range("N140").Value = -70988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 7680.3076
$ws.Range("I52").Value = 6863.5454
$ws.Range("J52").Value = 12172.5
$ws.Range("K52").Value = 6863.5454
$ws.Range("L52").Value = 12172.5
$ws.Range("M52").Value = -6637.5454
$ws.Range("N52").Value = -12624.5
$ws.Range("H62").Value = 6461
$ws.Range("I62").Value = 6744.1665
$ws.Range("J62").Value = 5894.6665
$ws.Range("K62").Value = 6744.1665
$ws.Range("L62").Value = 5894.6665
$ws.Range("M62").Value = -6120.1665
$ws.Range("N62").Value = -7142.6665
$ws.Range("H65").Value = 6461
$ws.Range("I65").Value = 6744.1665
$ws.Range("J65").Value = 5894.6665
$ws.Range("K65").Value = 33720.8325
$ws.Range("L65").Value = 29473.3325
$ws.Range("M65").Value = -30600.8325
$ws.Range("N65").Value = -35713.3325
$ws.Range("H81").Value = 3998.6667
$ws.Range("I81").Value = 4314.3335
$ws.Range("J81").Value = 3367.3333
$ws.Range("K81").Value = 8628.666999999999
$ws.Range("L81").Value = 6734.6666
$ws.Range("M81").Value = -7567.666999999999
$ws.Range("N81").Value = -8856.6666
$ws.Range("H84").Value = 3998.6667
$ws.Range("I84").Value = 4314.3335
$ws.Range("J84").Value = 3367.3333
$ws.Range("K84").Value = 43143.335
$ws.Range("L84").Value = 33673.333
$ws.Range("M84").Value = -37839.335
$ws.Range("N84").Value = -44281.333
$ws.Range("H94").Value = 105321.25
$ws.Range("J94").Value = 105321.25
$ws.Range("L94").Value = 105321.25
$ws.Range("N94").Value = -107123.25
$ws.Range("H95").Value = 125448
$ws.Range("J95").Value = 125448
$ws.Range("L95").Value = 125448
$ws.Range("N95").Value = -130940
$ws.Range("H96").Value = 7587.8887
$ws.Range("I96").Value = 1798.3334
$ws.Range("J96").Value = 10482.667
$ws.Range("K96").Value = 1798.3334
$ws.Range("L96").Value = 10482.667
$ws.Range("M96").Value = -425.3334
$ws.Range("N96").Value = -13228.667
$ws.Range("H100").Value = 1124.25
$ws.Range("I100").Value = 749.5
$ws.Range("J100").Value = 1499
$ws.Range("K100").Value = 1499
$ws.Range("L100").Value = 2998
$ws.Range("M100").Value = -958
$ws.Range("N100").Value = -4080
$ws.Range("H122").Value = 3336.6667
$ws.Range("J122").Value = 3436.8572
$ws.Range("L122").Value = 10310.5716
$ws.Range("N122").Value = -15210.5716
$ws.Range("H132").Value = 3892
$ws.Range("I132").Value = 2925.4119
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 8776.235700000001
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -6246.235700000001
$ws.Range("N132").Value = -29060
$ws.Range("H136").Value = 7377.385
$ws.Range("I136").Value = 14903.667
$ws.Range("J136").Value = 5119.5
$ws.Range("K136").Value = 44711.001
$ws.Range("L136").Value = 15358.5
$ws.Range("M136").Value = -42161.001
$ws.Range("N136").Value = -20458.5
